$wb = $excel.ActiveWorkbook

# --- Rename & reposition Sheet2 -> "硬件和操作系统" (moved to the front) ---
$tmp = $wb.Worksheets.Item("Sheet2")
$tmp.Name = "硬件和操作系统"
$tmp.Move($wb.Worksheets.Item(1))

# --- Sheet3 moves ahead of Sheet1 (order becomes 硬件和操作系统, Sheet3, Sheet1) ---
$tmp3 = $wb.Worksheets.Item("Sheet3")
$tmp3.Move($null, $wb.Worksheets.Item("硬件和操作系统"))

# Re-acquire the worksheet reference after Move (handle can go stale across Move)
$ws = $wb.Worksheets.Item("硬件和操作系统")

# --- Header row ---
$ws.Range("A1").Value = "服务器"
$ws.Range("B1").Value = "实例类型"
$ws.Range("C1").Value = "CPU(Core)"
$ws.Range("D1").Value = "内存（GB）"
$ws.Range("E1").Value = "磁盘分区-系统盘"
$ws.Range("F1").Value = "磁盘分区-数据盘"
$ws.Range("G1").Value = "挂载备份"
$ws.Range("H1").Value = "操作系统"

# --- Row 2: postgresql主库 ---
$ws.Range("A2").Value = "postgresql主库"
$ws.Range("B2").Value = "祼金属"
$ws.Range("C2").Value = 88
$ws.Range("D2").Value = 756
$ws.Range("E2").Value = "1T"
$ws.Range("F2").Value = "10T"
$ws.Range("G2").Value = "无"
$ws.Range("H2").Value = "CentOS release6.10"

# --- Row 3: postgresql备库 ---
$ws.Range("A3").Value = "postgresql备库"
$ws.Range("B3").Value = "祼金属"
$ws.Range("C3").Value = 88
$ws.Range("D3").Value = 756
$ws.Range("E3").Value = "1T"
$ws.Range("F3").Value = "10T"
$ws.Range("G3").Value = "35T"
$ws.Range("H3").Value = "CentOS release6.10"

# --- Row 4: postgresql异步备库 ---
$ws.Range("A4").Value = "postgresql异步备库"
$ws.Range("B4").Value = "云主机"
$ws.Range("G4").Value = "无"

# --- Column widths ---
$ws.Range("A1").EntireColumn.ColumnWidth = 18.5454545454545
$ws.Range("C1").EntireColumn.ColumnWidth = 9.72727272727273
$ws.Range("D1").EntireColumn.ColumnWidth = 9.81818181818182
$ws.Range("E1").EntireColumn.ColumnWidth = 15.8181818181818
$ws.Range("F1").EntireColumn.ColumnWidth = 16.9090909090909
$ws.Range("G1").EntireColumn.ColumnWidth = 16.9090909090909
$ws.Range("H1").EntireColumn.ColumnWidth = 19.5454545454545

# --- Borders around the full data range ---
$ws.Range("A1:H4").Borders.LineStyle = 1
$ws.Range("A1:H4").Borders.Weight = 2

# --- Selection like the saved file ---
$ws.Range("B9").Select()

$wb.Windows.Item(1).WindowState = -4143
